$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "height" column (E) values for rows 2-50 from centimeters to
# meters by dividing each existing value by 100.
for ($row = 2; $row -le 50; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # column E
    $cell.Value2 = $cell.Value2 / 100
}

# Update the selected/active cell shown in the sheet view.
$ws.Range("J8").Select()
